# "mise à jour des activités"
# Fill in the computed conversion result for row 2 of the Tableau1 table
# (24 GBP -> BEF), stored as text because column D ("resultat") uses the
# Text ("@") number format.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("D2").Value = "0,55008"
